$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "63.636.47"
$ws.Range("E2").Value = "  -1.05%  "
$ws.Range("D3").Value = "3.423.38"
$ws.Range("E3").Value = "  -2.26%  "
$ws.Range("E4").Value = "  +0.00%  "
$ws.Range("D5").Value = "'579.26"
$ws.Range("E5").Value = "  -1.97%  "
$ws.Range("D6").Value = "'129.25"
$ws.Range("E6").Value = "  -3.71%  "
$ws.Range("E7").Value = "  +0.05%  "
$ws.Range("D8").Value = "'0.480"
$ws.Range("E8").Value = "  -1.56%  "
$ws.Range("D9").Value = "'7.61"
$ws.Range("E9").Value = "  +3.89%  "
$ws.Range("D10").Value = "'0.123"
$ws.Range("E10").Value = "  -0.55%  "
$ws.Range("D11").Value = "'0.382"
$ws.Range("E11").Value = "  -1.09%  "
$ws.Range("D12").Value = "4.007.14"
$ws.Range("E12").Value = "  -2.21%  "
$ws.Range("E13").Value = "  -0.37%  "
$ws.Range("E14").Value = "  -2.57%  "
$ws.Range("D15").Value = "3.424.36"
$ws.Range("E15").Value = "  -2.21%  "
$ws.Range("D16").Value = "63.669.33"
$ws.Range("E16").Value = "  -1.00%  "
$ws.Range("D17").Value = "'25.14"
$ws.Range("E17").Value = "  -1.98%  "
$ws.Range("D18").Value = "'9.84"
$ws.Range("E18").Value = "  -0.11%  "
$ws.Range("D19").Value = "'5.64"
$ws.Range("E19").Value = "  -1.97%  "
$ws.Range("D20").Value = "'13.29"
$ws.Range("E20").Value = "  -1.71%  "
$ws.Range("D21").Value = "'383.10"
$ws.Range("E21").Value = "  -2.62%  "
$ws.Range("D22").Value = "'0.563"
$ws.Range("E22").Value = "  -1.63%  "
$ws.Range("D23").Value = "3.561.90"
$ws.Range("E23").Value = "  -2.19%  "
$ws.Range("E24").Value = "  -0.94%  "
$ws.Range("E25").Value = "  +0.14%  "
$ws.Range("E26").Value = "  -5.03%  "
$ws.Range("E27").Value = "  +0.16%  "
$ws.Range("E28").Value = "  -3.01%  "
$ws.Range("D29").Value = "'7.03"
$ws.Range("E29").Value = "  -4.81%  "
$ws.Range("D30").Value = "'7.92"
$ws.Range("E30").Value = "  -3.86%  "
$ws.Range("E31").Value = "  -0.87%  "
$ws.Range("D32").Value = "'1.41"
$ws.Range("E32").Value = "  -4.46%  "
$ws.Range("D33").Value = "3.454.70"
$ws.Range("E34").Value = "  -0.05%  "
$ws.Range("D35").Value = "'22.65"
$ws.Range("E35").Value = "  -3.52%  "
$ws.Range("D36").Value = "'5.14"
$ws.Range("E36").Value = "  +0.19%  "
$ws.Range("D37").Value = "'6.73"
$ws.Range("E37").Value = "  -2.12%  "
$ws.Range("D38").Value = "'164.11"
$ws.Range("E38").Value = "  -1.92%  "
$ws.Range("E39").Value = "  -2.58%  "
$ws.Range("D40").Value = "'0.0770"
$ws.Range("E40").Value = "  -1.28%  "
$ws.Range("E41").Value = "  -3.34%  "
$ws.Range("E42").Value = "  -0.01%  "
$ws.Range("D43").Value = "'41.37"
$ws.Range("E43").Value = "  -1.04%  "
$ws.Range("D44").Value = "'4.33"
$ws.Range("E44").Value = "  -1.63%  "
$ws.Range("E45").Value = "  -3.71%  "
$ws.Range("D46").Value = "'23.39"
$ws.Range("E46").Value = "  -7.01%  "
$ws.Range("E47").Value = "  -6.07%  "
$ws.Range("D48").Value = "'6.70"
$ws.Range("E48").Value = "  -0.78%  "
$ws.Range("E49").Value = "  -1.02%  "
$ws.Range("D50").Value = "2.286.69"
$ws.Range("E50").Value = "  -3.58%  "
$ws.Range("D51").Value = "'0.0252"
$ws.Range("E51").Value = "  -2.34%  "
